$wb = $excel.ActiveWorkbook

# ---- Sheet: Summary ----
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.6853932584269663
$ws1.Range("C2").Value = 0.6219211822660099
$ws1.Range("D2").Value = 0.9456928838951311
$ws1.Range("E2").Value = 0.75037147102526
$ws1.Range("F2").Value = 0.8565128900949797
$ws1.Range("G2").Value = 0.927128936590877
$ws1.Range("H2").Value = 0.7730996366900924
$ws1.Range("I2").Value = 505
$ws1.Range("J2").Value = 307
$ws1.Range("K2").Value = 227
$ws1.Range("L2").Value = 29

# ---- Sheet: Classification Report ----
$ws2 = $wb.Worksheets.Item("Classification Report")

$ws2.Range("B2").Value = 0.88671875
$ws2.Range("C2").Value = 0.4250936329588015
$ws2.Range("D2").Value = 0.5746835443037974

$ws2.Range("B3").Value = 0.6219211822660099
$ws2.Range("C3").Value = 0.9456928838951311
$ws2.Range("D3").Value = 0.75037147102526

$ws2.Range("B4").Value = 0.6853932584269663
$ws2.Range("C4").Value = 0.6853932584269663
$ws2.Range("D4").Value = 0.6853932584269663
$ws2.Range("E4").Value = 0.6853932584269663

$ws2.Range("B5").Value = 0.7543199661330049
$ws2.Range("C5").Value = 0.6853932584269663
$ws2.Range("D5").Value = 0.6625275076645287

$ws2.Range("B6").Value = 0.754319966133005
$ws2.Range("C6").Value = 0.6853932584269663
$ws2.Range("D6").Value = 0.6625275076645287

# ---- Sheet: Confusion Matrix ----
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

$ws3.Range("B2").Value = 227
$ws3.Range("C2").Value = 307

$ws3.Range("B3").Value = 29
$ws3.Range("C3").Value = 505
